# TDD: test1<Science>EntityEffect (for protractor, wheel, and tablet) pass
#
# Column D ("Tested?") on the Effects sheet gets filled in for the rows that
# were previously blank: "Y" for rows that were actually exercised by the
# test, "-" for rows that were not applicable / not exercised.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$testedYes = @(9, 10, 11, 12, 13, 15, 17, 19, 21, 23, 25)
$testedNo  = @(14, 16, 18, 20, 22, 24, 26, 27, 28)

foreach ($row in $testedYes) {
    $ws.Range("D$row").Value = "Y"
}

foreach ($row in $testedNo) {
    $ws.Range("D$row").Value = "-"
}

# Leave the selection where the author's session ended up.
$ws.Range("E34").Select() | Out-Null
